# Fix several data-entry errors in the "Checklist" sheet (EEPROM Container
# Review Checklist). Rows 18-24 get re-populated with corrected values and
# the stray duplicate row 25 is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checklist")

# --- Row 12: fix description text for Id_Config ---------------------------
$ws.Range("O12").Value = "description=Predefined data item for storing the configuration ID"

# --- Helper: clear a full data row (columns A:S) before rewriting it ------
function Clear-Row($row) {
    $ws.Range("A" + $row + ":S" + $row).ClearContents()
}

# --- Row 18: now TestData_08 ------------------------------------------------
Clear-Row 18
$ws.Range("A18").Value = "TestData_08"
$ws.Range("B18").Value = 12346
$ws.Range("C18").Value = 87
$ws.Range("E18").Value = "X"
$ws.Range("I18").Value = "X"
$ws.Range("M18").Value = "ee_range"
$ws.Range("O18").Value = "HOLA"
$ws.Range("P18").Value = 7
$ws.Range("Q18").Value = "JULIO"
$ws.Range("R18").Value = 8
$ws.Range("S18").Value = 5

# --- Row 19: now TestData_09 ------------------------------------------------
Clear-Row 19
$ws.Range("A19").Value = "TestData_09"
$ws.Range("B19").Value = 12347
$ws.Range("C19").Value = 4
$ws.Range("E19").Value = "X"
$ws.Range("F19").Value = "X"
$ws.Range("G19").Value = "X"
$ws.Range("M19").Value = "ee_range"
$ws.Range("N19").Value = "BB96"
$ws.Range("O19").Value = "ADIOS"
$ws.Range("P19").Value = 6
$ws.Range("Q19").Value = "RUBEN"
$ws.Range("R19").Value = 7
$ws.Range("S19").Value = 3

# --- Row 20: now TestData_10 ------------------------------------------------
Clear-Row 20
$ws.Range("A20").Value = "TestData_10"
$ws.Range("B20").Value = 12348
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = "X"
$ws.Range("E20").Value = "X"
$ws.Range("H20").Value = "X"
$ws.Range("I20").Value = "X"
$ws.Range("M20").Value = "ee_range"
$ws.Range("N20").Value = "BB75"
$ws.Range("O20").Value = ":("
$ws.Range("P20").Value = 5
$ws.Range("Q20").Value = "JULIO Y RUBEN"
$ws.Range("R20").Value = 6
$ws.Range("S20").Value = 6

# --- Row 21: now TestData_11 ------------------------------------------------
Clear-Row 21
$ws.Range("A21").Value = "TestData_11"
$ws.Range("B21").Value = 12349
$ws.Range("C21").Value = 5
$ws.Range("D21").Value = "X"
$ws.Range("F21").Value = "X"
$ws.Range("G21").Value = "X"
$ws.Range("H21").Value = "X"
$ws.Range("I21").Value = "X"
$ws.Range("M21").Value = "datablock"
$ws.Range("N21").Value = "BB89"
$ws.Range("P21").Value = 10
$ws.Range("Q21").Value = "JEJ"
$ws.Range("R21").Value = 5
$ws.Range("S21").Value = 9

# --- Row 22: now DUMMY_TestModuleCnt ---------------------------------------
Clear-Row 22
$ws.Range("A22").Value = "DUMMY_TestModuleCnt"
$ws.Range("B22").Value = 31416
$ws.Range("I22").Value = "X"
$ws.Range("M22").Value = "ee_erase"
$ws.Range("O22").Value = "description=- Component: DUMMY`n- REPROG info: use case REPROG must not be set (data must not be changed after reprogramming)!`nSometimes the description is longer than two rows.`nOther times, there are more than three.`nIn this case, it is one more than four. And could be more."

# --- Row 23: now ASDFClockTower --------------------------------------------
Clear-Row 23
$ws.Range("A23").Value = "ASDFClockTower"
$ws.Range("B23").Value = 111255
$ws.Range("G23").Value = "X"
$ws.Range("M23").Value = "ee_datablock"
$ws.Range("O23").Value = "description=- Component: ASDF`n- REPROG info: use case REPROG must not be set.`n- REPROG info: In certain cases there are two comments of this type.`nThere are also strings up to 160 characters per row, only on description fields and usually is not only one row. Like this example but a little bit longer."

# --- Row 24: now TestData_07 (brand-new corrected content) -----------------
Clear-Row 24
$ws.Range("A24").Value = "TestData_07"
$ws.Range("B24").Value = 123456
$ws.Range("C24").Value = 91
$ws.Range("F24").Value = "X"
$ws.Range("G24").Value = "X"
$ws.Range("H24").Value = "X"
$ws.Range("L24").Value = "X"
$ws.Range("M24").Value = "ee_range"
$ws.Range("S24").Value = 9

# --- Row 25: removed (its content was merged up into row 23/24) -----------
Clear-Row 25
